# Update the yearly exam statistics figures on Sheet1.
# This mirrors an "Add files via upload" re-export where several
# login-count (column C) and entry-count (column D) totals were bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 86
$ws.Range("C7").Value = 75
$ws.Range("C8").Value = 47
$ws.Range("D8").Value = 31
$ws.Range("D9").Value = 44
$ws.Range("D11").Value = 55
$ws.Range("C12").Value = 38
$ws.Range("C14").Value = 96
$ws.Range("D14").Value = 85
$ws.Range("D17").Value = 29
$ws.Range("C18").Value = 46
$ws.Range("C19").Value = 49
$ws.Range("D19").Value = 39
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 52
$ws.Range("D22").Value = 27
$ws.Range("C23").Value = 26
$ws.Range("D23").Value = 21
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = 32
$ws.Range("D30").Value = 45
$ws.Range("C32").Value = 45
$ws.Range("D32").Value = 37
$ws.Range("D33").Value = 37
$ws.Range("C36").Value = 40
$ws.Range("D36").Value = 26
$ws.Range("C40").Value = 82
$ws.Range("D40").Value = 72
$ws.Range("C41").Value = 62
$ws.Range("D41").Value = 51
$ws.Range("C43").Value = 63
$ws.Range("D43").Value = 54
$ws.Range("C44").Value = 69
$ws.Range("D44").Value = 58
$ws.Range("C48").Value = 52
$ws.Range("D48").Value = 35
$ws.Range("C52").Value = 52
$ws.Range("D54").Value = 35
$ws.Range("C57").Value = 70
$ws.Range("C63").Value = 106
$ws.Range("D63").Value = 98
$ws.Range("C66").Value = 65
$ws.Range("C68").Value = 60
$ws.Range("D68").Value = 47
$ws.Range("C70").Value = 53
$ws.Range("C73").Value = 59
$ws.Range("C77").Value = 134
$ws.Range("D77").Value = 134
$ws.Range("C78").Value = 80
$ws.Range("D78").Value = 63
$ws.Range("C81").Value = 75
$ws.Range("C82").Value = 29
$ws.Range("C83").Value = 71
$ws.Range("D83").Value = 63
$ws.Range("D84").Value = 127
$ws.Range("C92").Value = 249
$ws.Range("D92").Value = 183
$ws.Range("C93").Value = 5505
$ws.Range("D93").Value = 4438
